# Apply the change described by the diff:
# Insert a new task row ("validation de numero d'emplyee dans le login ")
# just above the totals rows, shifting the totals down by one row, and
# updating the SUM/aggregate formulas to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlNone = -4142

# Insert a new row at 22 (pushes old row22/23 -> 23/24), inheriting the
# formatting of the row above (row 21).
$ws.Rows.Item(22).Insert()

# --- Fill in the new task row (row 22) ---
$ws.Range("B22").Value = "validation de numero d'emplyee dans le login "
$ws.Range("C22").Value = 1
$ws.Range("D22").Formula = "=C22*60"

# The new row does not carry the heavy border used by the task rows above
# it (only a light wrap-text formatting remains on D22), so drop the
# borders that were inherited from row 21 and restore plain formatting.
$ws.Range("B22:F22").Borders.LineStyle = $xlNone
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").WrapText = $true
$ws.Range("E22:F22").IndentLevel = 0

# Make sure the row keeps the same visual height as the rest of the sheet.
$ws.Rows.Item(22).RowHeight = 43.95

# --- Fix up the totals rows (now shifted to 23 and 24) ---
$ws.Range("C23").Formula = "=SUM(C3:C22)"
$ws.Range("D23").Formula = "=SUM(D3:D22)"
$ws.Range("D24").Formula = "=D23/60"

# --- Cosmetic view state (best effort) ---
$ws.PageSetup.Orientation = 1
$ws.Range("H24").Select()

Write-Host "Row inserted and totals updated"
